$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Object__to_from_node")

# Remove the fom_cost / vom_cost columns (G:H) entirely - shifts nothing left of G, just removes them
$ws.Range("G1:H28").Delete()

# Rows 20-28 are no longer part of the table; clear their contents (A:F)
$ws.Range("A20:F28").ClearContents()

# Rewrite rows 2-19 (A:F) with the updated relationship data
$ws.Cells.Item(2,1).Value = "unit__to_node"
$ws.Cells.Item(2,2).Value = "unit"
$ws.Cells.Item(2,3).Value = "Solar_Plant_Kasso"
$ws.Cells.Item(2,4).Value = "Power_Kasso"
$ws.Cells.Item(2,5).Value = "unit_capacity"
$ws.Cells.Item(2,6).Value = 304
$ws.Cells.Item(3,1).Value = "unit__from_node"
$ws.Cells.Item(3,2).Value = "unit"
$ws.Cells.Item(3,3).Value = "Electrolyzer"
$ws.Cells.Item(3,4).Value = "Power_Kasso"
$ws.Cells.Item(3,5).Value = "unit_capacity"
$ws.Cells.Item(3,6).Value = 52
$ws.Cells.Item(4,1).Value = "unit__from_node"
$ws.Cells.Item(4,2).Value = "unit"
$ws.Cells.Item(4,3).Value = "Electrolyzer"
$ws.Cells.Item(4,4).Value = "Power_Kasso"
$ws.Cells.Item(4,5).Value = "vom_cost"
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(5,1).Value = "unit__to_node"
$ws.Cells.Item(5,2).Value = "unit"
$ws.Cells.Item(5,3).Value = "CO2_Vaporizer"
$ws.Cells.Item(5,4).Value = "Vaporized_Carbon_Dioxide"
$ws.Cells.Item(5,5).Value = "unit_capacity"
$ws.Cells.Item(5,6).Value = 100
$ws.Cells.Item(6,1).Value = "unit__to_node"
$ws.Cells.Item(6,2).Value = "unit"
$ws.Cells.Item(6,3).Value = "Destilation_Tower"
$ws.Cells.Item(6,4).Value = "E-Methanol_Kasso"
$ws.Cells.Item(6,5).Value = "unit_capacity"
$ws.Cells.Item(6,6).Value = 52
$ws.Cells.Item(7,1).Value = "unit__to_node"
$ws.Cells.Item(7,2).Value = "unit"
$ws.Cells.Item(7,3).Value = "Methanol_Reactor"
$ws.Cells.Item(7,4).Value = "Raw_Methanol"
$ws.Cells.Item(7,5).Value = "unit_capacity"
$ws.Cells.Item(7,6).Value = 100
$ws.Cells.Item(8,1).Value = "unit__to_node"
$ws.Cells.Item(8,2).Value = "unit"
$ws.Cells.Item(8,3).Value = "Methanol_Reactor"
$ws.Cells.Item(8,4).Value = "Waste_Heat"
$ws.Cells.Item(8,5).Value = "unit_capacity"
$ws.Cells.Item(8,6).Value = 100
$ws.Cells.Item(9,1).Value = "connection__from_node"
$ws.Cells.Item(9,2).Value = "connection"
$ws.Cells.Item(9,3).Value = "power_line_Wholesale_Kasso"
$ws.Cells.Item(9,4).Value = "Power_Wholesale"
$ws.Cells.Item(9,5).Value = "connection_capacity"
$ws.Cells.Item(9,6).Value = 1000
$ws.Cells.Item(10,1).Value = "connection__to_node"
$ws.Cells.Item(10,2).Value = "connection"
$ws.Cells.Item(10,3).Value = "power_line_Wholesale_Kasso"
$ws.Cells.Item(10,4).Value = "Power_Kasso"
$ws.Cells.Item(10,5).Value = "connection_capacity"
$ws.Cells.Item(10,6).Value = 1000
$ws.Cells.Item(11,1).Value = "connection__from_node"
$ws.Cells.Item(11,2).Value = "connection"
$ws.Cells.Item(11,3).Value = "power_line_Wholesale_Kasso"
$ws.Cells.Item(11,4).Value = "Power_Kasso"
$ws.Cells.Item(11,5).Value = "connection_capacity"
$ws.Cells.Item(11,6).Value = 1000
$ws.Cells.Item(12,1).Value = "connection__to_node"
$ws.Cells.Item(12,2).Value = "connection"
$ws.Cells.Item(12,3).Value = "power_line_Wholesale_Kasso"
$ws.Cells.Item(12,4).Value = "Power_Wholesale"
$ws.Cells.Item(12,5).Value = "connection_capacity"
$ws.Cells.Item(12,6).Value = 1000
$ws.Cells.Item(13,1).Value = "connection__from_node"
$ws.Cells.Item(13,2).Value = "connection"
$ws.Cells.Item(13,3).Value = "pipeline_storage_hydrogen"
$ws.Cells.Item(13,4).Value = "Hydrogen_Kasso"
$ws.Cells.Item(13,5).Value = "vom_cost"
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(14,1).Value = "connection__to_node"
$ws.Cells.Item(14,2).Value = "connection"
$ws.Cells.Item(14,3).Value = "pipeline_storage_hydrogen"
$ws.Cells.Item(14,4).Value = "Hydrogen_storage_Kasso"
$ws.Cells.Item(14,5).Value = "connection_capacity"
$ws.Cells.Item(14,6).Value = 1000
$ws.Cells.Item(15,1).Value = "connection__from_node"
$ws.Cells.Item(15,2).Value = "connection"
$ws.Cells.Item(15,3).Value = "pipeline_storage_hydrogen"
$ws.Cells.Item(15,4).Value = "Hydrogen_storage_Kasso"
$ws.Cells.Item(15,5).Value = "connection_capacity"
$ws.Cells.Item(15,6).Value = 1000
$ws.Cells.Item(16,1).Value = "connection__to_node"
$ws.Cells.Item(16,2).Value = "connection"
$ws.Cells.Item(16,3).Value = "pipeline_storage_e-methanol"
$ws.Cells.Item(16,4).Value = "E-Methanol_storage_Kasso"
$ws.Cells.Item(16,5).Value = "connection_capacity"
$ws.Cells.Item(16,6).Value = 1000
$ws.Cells.Item(17,1).Value = "connection__from_node"
$ws.Cells.Item(17,2).Value = "connection"
$ws.Cells.Item(17,3).Value = "pipeline_storage_e-methanol"
$ws.Cells.Item(17,4).Value = "E-Methanol_storage_Kasso"
$ws.Cells.Item(17,5).Value = "connection_capacity"
$ws.Cells.Item(17,6).Value = 1000
$ws.Cells.Item(18,1).Value = "connection__from_node"
$ws.Cells.Item(18,2).Value = "connection"
$ws.Cells.Item(18,3).Value = "pipeline_District_Heating"
$ws.Cells.Item(18,4).Value = "Waste_Heat"
$ws.Cells.Item(18,5).Value = "connection_capacity"
$ws.Cells.Item(18,6).Value = 1000
$ws.Cells.Item(19,1).Value = "connection__to_node"
$ws.Cells.Item(19,2).Value = "connection"
$ws.Cells.Item(19,3).Value = "pipeline_District_Heating"
$ws.Cells.Item(19,4).Value = "District_Heating"
$ws.Cells.Item(19,5).Value = "connection_capacity"
$ws.Cells.Item(19,6).Value = 1000
Write-Output "edit applied"
